$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used by existing text data cells (e.g. Coin/Price/Volume columns),
# so that re-assigning values does not introduce a new implicit "Text" number format style.
$dataStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'64.025.34"
$ws.Range("D2").Style = $dataStyle
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "'2.753.86"
$ws.Range("D3").Style = $dataStyle
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'577.80"
$ws.Range("D5").Style = $dataStyle
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'157.92"
$ws.Range("D6").Style = $dataStyle
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.607"
$ws.Range("D8").Style = $dataStyle
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").Value = "'5.85"
$ws.Range("D10").Style = $dataStyle
$ws.Range("E10").Value = "  -12.93%  "
$ws.Range("D11").Value = "'0.386"
$ws.Range("D11").Style = $dataStyle
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "'3.243.55"
$ws.Range("D13").Style = $dataStyle
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "'26.82"
$ws.Range("D14").Style = $dataStyle
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "'63.927.00"
$ws.Range("D15").Style = $dataStyle
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'2.755.83"
$ws.Range("D17").Style = $dataStyle
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'12.10"
$ws.Range("D18").Style = $dataStyle
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'4.88"
$ws.Range("D19").Style = $dataStyle
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'359.72"
$ws.Range("D20").Style = $dataStyle
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("D21").Style = $dataStyle
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "'0.552"
$ws.Range("D22").Style = $dataStyle
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = $dataStyle
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'66.21"
$ws.Range("D24").Style = $dataStyle
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "'8.48"
$ws.Range("D26").Style = $dataStyle
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'0.0₃0931"
$ws.Range("D28").Style = $dataStyle
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "'1.96"
$ws.Range("D29").Style = $dataStyle
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "'7.03"
$ws.Range("D30").Style = $dataStyle
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").Value = "'169.91"
$ws.Range("D32").Style = $dataStyle
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("D33").Value = "'20.36"
$ws.Range("D33").Style = $dataStyle
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'4.93"
$ws.Range("D34").Style = $dataStyle
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = $dataStyle
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'0.997"
$ws.Range("D38").Style = $dataStyle
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'6.19"
$ws.Range("D39").Style = $dataStyle
$ws.Range("E39").Value = "  +11.56%  "
$ws.Range("D40").Value = "'4.17"
$ws.Range("D40").Style = $dataStyle
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "'331.00"
$ws.Range("D41").Style = $dataStyle
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("D42").Value = "'39.19"
$ws.Range("D42").Style = $dataStyle
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "'21.76"
$ws.Range("D43").Style = $dataStyle
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'0.0593"
$ws.Range("D44").Style = $dataStyle
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "'21.76"
$ws.Range("D45").Style = $dataStyle
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.635"
$ws.Range("D46").Style = $dataStyle
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0256"
$ws.Range("D47").Style = $dataStyle
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'136.19"
$ws.Range("D48").Style = $dataStyle
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +0.71%  "
